$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 206.55556
$ws.Range("I9").Value = 194.8
$ws.Range("J9").Value = 221.25
$ws.Range("K9").Value = 194.8
$ws.Range("L9").Value = 221.25
$ws.Range("M9").Value = -25.80000000000001
$ws.Range("N9").Value = -559.25

$ws.Range("H19").Value = 3687
$ws.Range("I19").Value = 6997.2
$ws.Range("J19").Value = 928.5
$ws.Range("K19").Value = 6997.2
$ws.Range("L19").Value = 928.5
$ws.Range("M19").Value = -6822.2
$ws.Range("N19").Value = -1278.5

$ws.Range("H74").Value = 8683.3125
$ws.Range("I74").Value = 6757.143
$ws.Range("K74").Value = 6757.143
$ws.Range("M74").Value = -5821.143

$ws.Range("H77").Value = 8683.3125
$ws.Range("I77").Value = 6757.143
$ws.Range("K77").Value = 33785.715
$ws.Range("M77").Value = -29105.715

$ws.Range("H107").Value = 2220.5312
$ws.Range("J107").Value = 3247.625
$ws.Range("L107").Value = 3247.625
$ws.Range("N107").Value = -7087.625

$ws.Range("H113").Value = 6350.773
$ws.Range("I113").Value = 3412.5
$ws.Range("K113").Value = 3412.5
$ws.Range("M113").Value = -158.5

$ws.Range("H116").Value = 16000
$ws.Range("J116").Value = 15400
$ws.Range("L116").Value = 15400
$ws.Range("N116").Value = -22284

$ws.Range("H138").Value = 3133.5
$ws.Range("J138").Value = 4170.6816
$ws.Range("L138").Value = 12512.0448
$ws.Range("N138").Value = -22792.0448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 43480696
$ws.Range("I45").Value = 55556850
$ws.Range("J45").Value = 6540.4
$ws.Range("K45").Value = 55556850
$ws.Range("L45").Value = 6540.4
$ws.Range("M45").Value = -55556473
$ws.Range("N45").Value = -7294.4

$ws.Range("H61").Value = 7789.4614
$ws.Range("I61").Value = 6788.625
$ws.Range("K61").Value = 6788.625
$ws.Range("M61").Value = -6576.625

$ws.Range("H74").Value = 20835586
$ws.Range("I74").Value = 23811852
$ws.Range("K74").Value = 23811852
$ws.Range("M74").Value = -23810978

$ws.Range("H77").Value = 20835586
$ws.Range("I77").Value = 23811852
$ws.Range("K77").Value = 119059260
$ws.Range("M77").Value = -119054892

$ws.Range("H122").Value = 1335.591
$ws.Range("I122").Value = 623.375
$ws.Range("K122").Value = 1870.125
$ws.Range("M122").Value = 579.875

$ws.Range("H136").Value = 7789.4614
$ws.Range("I136").Value = 6788.625
$ws.Range("K136").Value = 20365.875
$ws.Range("M136").Value = -17815.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2918.2222
$ws.Range("I86").Value = 2918.2222
$ws.Range("K86").Value = 2918.2222
$ws.Range("M86").Value = -1795.2222

$ws.Range("H89").Value = 2918.2222
$ws.Range("I89").Value = 2918.2222
$ws.Range("K89").Value = 14591.111
$ws.Range("M89").Value = -8975.111000000001

$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1549
$ws.Range("N94").ClearContents()

$ws.Range("H107").Value = 1223.04
$ws.Range("I107").Value = 1187.1052
$ws.Range("J107").Value = 1336.8334
$ws.Range("K107").Value = 1187.1052
$ws.Range("L107").Value = 1336.8334
$ws.Range("M107").Value = 732.8948
$ws.Range("N107").Value = -5176.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1808.4546
$ws.Range("I16").Value = 1208
$ws.Range("K16").Value = 1208
$ws.Range("M16").Value = -921

$ws.Range("H107").Value = 904.1923
$ws.Range("I107").Value = 805.65
$ws.Range("J107").Value = 1232.6666
$ws.Range("K107").Value = 805.65
$ws.Range("L107").Value = 1232.6666
$ws.Range("M107").Value = 1114.35
$ws.Range("N107").Value = -5072.6666

$ws.Range("H113").Value = 1808.4546
$ws.Range("I113").Value = 1208
$ws.Range("K113").Value = 1208
$ws.Range("M113").Value = 962

$ws.Range("H132").Value = 8331.333000000001
$ws.Range("I132").Value = 8662.5
$ws.Range("K132").Value = 25987.5
$ws.Range("M132").Value = -23457.5

$ws.Range("H134").Value = 3418.9443
$ws.Range("I134").Value = 2101.1667
$ws.Range("K134").Value = 6303.500100000001
$ws.Range("M134").Value = -3768.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3539.889
$ws.Range("J68").Value = 3669.875
$ws.Range("L68").Value = 11009.625
$ws.Range("N68").Value = -12631.625

$ws.Range("H71").Value = 3539.889
$ws.Range("J71").Value = 3669.875
$ws.Range("L71").Value = 33028.875
$ws.Range("N71").Value = -41140.875

$ws.Range("H132").Value = 5150.25
$ws.Range("I132").Value = 3733.3333
$ws.Range("K132").Value = 33599.9997
$ws.Range("M132").Value = -31069.9997

$ws.Range("H134").Value = 5431.357
$ws.Range("I134").Value = 829.7143
$ws.Range("J134").Value = 10033
$ws.Range("K134").Value = 2489.1429
$ws.Range("L134").Value = 30099
$ws.Range("M134").Value = 2580.8571
$ws.Range("N134").Value = -40239

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 29668.666
$ws.Range("I21").Value = 29668.666
$ws.Range("K21").Value = 29668.666
$ws.Range("M21").Value = -29495.666

$ws.Range("H30").Value = 29668.666
$ws.Range("I30").Value = 29668.666
$ws.Range("K30").Value = 29668.666
$ws.Range("M30").Value = -29563.666

$ws.Range("H43").Value = 3716.4443
$ws.Range("I43").Value = 350
$ws.Range("J43").Value = 15499
$ws.Range("K43").Value = 350
$ws.Range("L43").Value = 15499
$ws.Range("M43").Value = -199
$ws.Range("N43").Value = -15801

$ws.Range("H98").Value = 20000
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 20000
$ws.Range("N98").Value = -25990

$ws.Range("H122").Value = 3324.8333
$ws.Range("I122").Value = 2918.25
$ws.Range("K122").Value = 8754.75
$ws.Range("M122").Value = -6304.75

$ws.Range("H126").Value = 4255.5
$ws.Range("I126").Value = 3343.5386
$ws.Range("K126").Value = 10030.6158
$ws.Range("M126").Value = -7560.6158

$ws.Range("H131").Value = 59998.5
$ws.Range("J131").Value = 59998.5
$ws.Range("L131").Value = 59998.5
$ws.Range("N131").Value = -70078.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9908.571
$ws.Range("I7").Value = 4819.6
$ws.Range("K7").Value = 4819.6
$ws.Range("M7").Value = -4707.6

$ws.Range("H22").Value = 7617
$ws.Range("I22").Value = 1766.6666
$ws.Range("J22").Value = 9810.875
$ws.Range("K22").Value = 1766.6666
$ws.Range("L22").Value = 9810.875
$ws.Range("M22").Value = -1471.6666
$ws.Range("N22").Value = -10400.875

$ws.Range("H27").Value = 7617
$ws.Range("I27").Value = 1766.6666
$ws.Range("J27").Value = 9810.875
$ws.Range("K27").Value = 1766.6666
$ws.Range("L27").Value = 9810.875
$ws.Range("M27").Value = -1659.6666
$ws.Range("N27").Value = -10024.875

$ws.Range("H46").Value = 2529.9
$ws.Range("I46").Value = 625.5
$ws.Range("K46").Value = 625.5
$ws.Range("M46").Value = -437.5

$ws.Range("H93").Value = 2892.3333
$ws.Range("I93").Value = 589
$ws.Range("J93").Value = 7499
$ws.Range("K93").Value = 589
$ws.Range("L93").Value = 7499
$ws.Range("M93").Value = 659
$ws.Range("N93").Value = -9995

$ws.Range("H126").Value = 9908.571
$ws.Range("I126").Value = 4819.6
$ws.Range("K126").Value = 14458.8
$ws.Range("M126").Value = -11988.8

$ws.Range("H132").Value = 4930.8
$ws.Range("I132").Value = 2518.261
$ws.Range("K132").Value = 7554.782999999999
$ws.Range("M132").Value = -5024.782999999999

$ws.Range("H136").Value = 8541.65
$ws.Range("J136").Value = 12131.308
$ws.Range("L136").Value = 36393.924
$ws.Range("N136").Value = -41493.924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H122").Value = 9056.038
$ws.Range("I122").Value = 2866.3333
$ws.Range("J122").Value = 12332.941
$ws.Range("K122").Value = 8598.999899999999
$ws.Range("L122").Value = 36998.823
$ws.Range("M122").Value = -6148.999899999999
$ws.Range("N122").Value = -41898.823

$ws.Range("H132").Value = 6397.7666
$ws.Range("I132").Value = 5664
$ws.Range("K132").Value = 16992
$ws.Range("M132").Value = -14462
